$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.06"
$ws.Range("E2").Value = "'0.23%"
$ws.Range("D3").Value = "'40.95"
$ws.Range("E3").Value = "'2.35%"
$ws.Range("D4").Value = "'5.102"
$ws.Range("E4").Value = "'-0.49%"
$ws.Range("D5").Value = "'0.07625"
$ws.Range("E5").Value = "'-1.41%"
$ws.Range("E6").Value = "'0.35%"
$ws.Range("D7").Value = "'1.606"
$ws.Range("E7").Value = "'-0.54%"
$ws.Range("D9").Value = "'0.9024"
$ws.Range("E9").Value = "'1.78%"
$ws.Range("D10").Value = "'0.1121"
$ws.Range("E10").Value = "'12.78%"
$ws.Range("D11").Value = "'0.1782"
$ws.Range("E11").Value = "'2.55%"
$ws.Range("D12").Value = "'0.09132"
$ws.Range("E12").Value = "'0.75%"
$ws.Range("D13").Value = "'0.04209"
$ws.Range("E13").Value = "'-5.32%"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.001260"
$ws.Range("E15").Value = "'0.23%"
$ws.Range("D16").Value = "'0.005714"
$ws.Range("E16").Value = "'-3.87%"
$ws.Range("E17").Value = "'-0.09%"
$ws.Range("E18").Value = "'0.65%"
$ws.Range("D19").Value = "'6.625"
$ws.Range("E19").Value = "'-6.88%"
$ws.Range("D20").Value = "'0.1366"
$ws.Range("E20").Value = "'1.25%"
$ws.Range("D22").Value = "'0.04069"
$ws.Range("E22").Value = "'-1.42%"
$ws.Range("D23").Value = "'0.001244"
$ws.Range("E23").Value = "'3.59%"
$ws.Range("D24").Value = "'0.004118"
$ws.Range("E24").Value = "'1.20%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'0.04%"
$ws.Range("E26").Value = "'-94.99%"
$ws.Range("D38").Value = "'0.02380"
$ws.Range("E38").Value = "'1.27%"
$ws.Range("D39").Value = "'0.05172"
$ws.Range("E39").Value = "'-0.62%"
$ws.Range("D40").Value = "'0.007783"
$ws.Range("E40").Value = "'-1.89%"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("E41").Value = "'-1.77%"
$ws.Range("D42").Value = "'0.007057"
$ws.Range("E42").Value = "'12.46%"
$ws.Range("D43").Value = "'0.001952"
$ws.Range("E43").Value = "'0.06%"
$ws.Range("D44").Value = "'0.007733"
$ws.Range("E44").Value = "'-11.58%"
$ws.Range("D45").Value = "'0.3080"
$ws.Range("E45").Value = "'-7.56%"
$ws.Range("D46").Value = "'0.00007017"
$ws.Range("E46").Value = "'7.00%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'0.05568"
$ws.Range("E48").Value = "'1,469.46%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.07%"
